$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit inserts one new product row ("كريم براشوت", price 45) into the
# price-list table right after "كريم ONE للبشره الحساسه" (currently row 93),
# pushing the rows below it (93..99) down to (94..100), and bumps the
# running total in the totals row by the new product's price (45).
#
# Observed target shape (diffed from the canonical OOXML):
#   - Columns B/C/D/E/F/G/H/I/J/K/L/M/N shift down by one row for the
#     affected range; the brand new row reuses whatever H/N text the
#     row-that-is-now-93 already had (consistent with "insert a duplicate
#     row, then edit the name + price"), and gets a new name + price.
#   - Column A (the sequential "م" index) is NOT shifted along with the
#     rest: it is simply `row - 3` for every data row, so the existing
#     values already in place stay put and only the brand new row at the
#     bottom of the shifted block needs a new value appended.
#   - Row heights are independent, per-row-index values (not tied to a
#     row's content), so only the rows whose role changed (old totals/
#     footer rows, now one row lower) get an explicit height.
#
# We avoid native Rows.Insert(): every data row already shares identical
# per-column cell styles, so a plain value shift (no structural insert)
# reproduces the desired result without perturbing styles.xml with
# extra/duplicate style entries the way a literal row insert would.
# ---------------------------------------------------------------------------

$firstRow = 93
$lastRow = 99
$lastCol = 14  # column N

# Capture the current (pre-edit) values for every column except A, for the
# rows that will shift down (93..99), before anything is overwritten.
$captured = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $captured["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Clear out the text cells whose shared-string entries must be regenerated
# (a shared string is garbage-collected once its last reference disappears),
# so that re-writing them further down, in the correct order, rebuilds the
# shared-string table in the same sequence the real edit produced: new
# product name first, then each shifted string in its original relative
# order.
$ws.Cells.Item(93, 2).Value2 = ""   # B93  ماء اكسجين 20
$ws.Cells.Item(94, 2).Value2 = ""   # B94  ماسك جلسات اطفال
$ws.Cells.Item(95, 2).Value2 = ""   # B95  معجون سيجنال 25 مل
$ws.Cells.Item(96, 2).Value2 = ""   # B96  معجون سيجنال عرض 50ملل
$ws.Cells.Item(97, 2).Value2 = ""   # B97  معجون كلوز اب الصغير
$ws.Cells.Item(97, 8).Value2 = ""   # H97  16:0
$ws.Cells.Item(99, 1).Value2 = ""   # A99  Wednesday, 7 January, 2026 7:22 PM
$ws.Cells.Item(99, 6).Value2 = ""   # F99  1/1
$ws.Cells.Item(99, 9).Value2 = ""   # I99  developed by : Abdelaziz Talaat

# Shift columns B..N down by one row (destination = source + 1), for every
# column except the two (B, H) whose text needs to be re-entered in a
# specific order afterwards to control shared-string placement.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $dst = $r + 1
    for ($c = 2; $c -le $lastCol; $c++) {
        if ($c -eq 2 -or $c -eq 8) {
            continue
        }
        $ws.Cells.Item($dst, $c).Value2 = $captured["$r,$c"]
    }
}

# Re-write the shifted text columns (B, H) strictly in ascending
# destination-row order, so new shared-string entries are appended in the
# sequence: كريم براشوت, ماء اكسجين 20, ماسك جلسات اطفال,
# معجون سيجنال 25 مل, معجون سيجنال عرض 50ملل, معجون كلوز اب الصغير, 16:0,
# <date>, 1/1, developed by : ...
$ws.Cells.Item(93, 2).Value2 = "كريم براشوت"
$ws.Cells.Item(94, 2).Value2 = $captured["93,2"]
$ws.Cells.Item(95, 2).Value2 = $captured["94,2"]
$ws.Cells.Item(96, 2).Value2 = $captured["95,2"]
$ws.Cells.Item(97, 2).Value2 = $captured["96,2"]
$ws.Cells.Item(98, 2).Value2 = $captured["97,2"]
$ws.Cells.Item(98, 8).Value2 = $captured["97,8"]
$ws.Cells.Item(100, 1).Value2 = $captured["99,1"]
$ws.Cells.Item(100, 6).Value2 = $captured["99,6"]
$ws.Cells.Item(100, 9).Value2 = $captured["99,9"]

# New row 93 ("كريم براشوت"): only the name (B) and price (L) are genuinely
# new data; H93/N93 are left untouched because they already hold the exact
# values the new row needs (the row that is being duplicated-in-place).
$ws.Cells.Item(93, 12).Value2 = 45

# Column A ("م" / sequential index): every other data row already carries
# its correct `row - 3` value, only the newly-appended row at the bottom of
# the shifted block needs one.
$ws.Cells.Item(98, 1).Value2 = 95

# Update the running total (row 99 now, previously row 98) by the new
# product's price.
$ws.Cells.Item(99, 11).Value2 = $captured["98,11"] + 45

# ---------------------------------------------------------------------------
# Merged cells: the merge ranges for the old totals row (98) and footer row
# (99) need to move down to (99) and (100) respectively, and the newly
# promoted data row (98) needs the same 3-way merge pattern
# (B:G / H:K / L:M) every other product row uses.
# ---------------------------------------------------------------------------
$ws.Range("K98:N98").UnMerge()
$ws.Range("A99:E99").UnMerge()
$ws.Range("F99:G99").UnMerge()
$ws.Range("I99:N99").UnMerge()

$ws.Range("B98:G98").Merge()
$ws.Range("H98:K98").Merge()
$ws.Range("L98:M98").Merge()
$ws.Range("K99:N99").Merge()
$ws.Range("A100:E100").Merge()
$ws.Range("F100:G100").Merge()
$ws.Range("I100:N100").Merge()

# Row heights: every other data row keeps its own original height (heights
# are per-row-index, not tied to content), so only the rows whose role
# changed need an explicit height.
$ws.Rows(99).RowHeight = 25.5   # now the totals row
$ws.Rows(100).RowHeight = 16.5  # now the footer row

Write-Host "Inserted new product row; total updated."
